# Apply updated price/volume values to the cryptos worksheet.
# Each D/E cell holds plain text (e.g. "278.42", "1.05%"), so the
# NumberFormat is forced to "@" (Text) before assigning the new
# value -- otherwise Excel COM auto-coerces numeric-looking strings
# into real numbers (and "1.05%" into a computed percentage).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.05%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.36%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.837"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.09%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06373"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.53%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.029"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.26%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.308"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.09%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8930"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.92%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1519"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.42%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05588"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "10.63%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07443"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.64%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02938"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.57%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08964"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.67%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.15%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006359"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.46%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006116"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.85%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.83%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.321"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.73%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.233"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.71%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1348"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.34%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.21%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04387"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.15%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.43%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004277"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "10.85%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001179"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.63%"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-8.28%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04029"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.10%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006727"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.76%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1408"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "19.65%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002078"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.45%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01116"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.19%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005541"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.30%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.561"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.01%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01847"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-19.58%"
